# Auto-generated Excel COM-interop edit script
# Applies numeric updates to leve-profit tables across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as described by the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4351.231
$ws.Range("J17").Value = 3585.4546
$ws.Range("L17").Value = 10756.3638
$ws.Range("N17").Value = -11092.3638
$ws.Range("H132").Value = 944.9216
$ws.Range("J132").Value = 2197.5
$ws.Range("L132").Value = 6592.5
$ws.Range("N132").Value = -11652.5
$ws.Range("H138").Value = 3156.366
$ws.Range("J138").Value = 2806.8572
$ws.Range("L138").Value = 8420.571599999999
$ws.Range("N138").Value = -18700.5716

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 15000
$ws.Range("J55").Value = 15000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15630
$ws.Range("H74").Value = 997.73914
$ws.Range("I74").Value = 753.9268
$ws.Range("J74").Value = 2997
$ws.Range("K74").Value = 753.9268
$ws.Range("L74").Value = 2997
$ws.Range("M74").Value = 120.0732
$ws.Range("N74").Value = -4745
$ws.Range("H77").Value = 997.73914
$ws.Range("I77").Value = 753.9268
$ws.Range("J77").Value = 2997
$ws.Range("K77").Value = 3769.634
$ws.Range("L77").Value = 14985
$ws.Range("M77").Value = 598.366
$ws.Range("N77").Value = -23721
$ws.Range("H86").Value = 25000
$ws.Range("J86").Value = 25000
$ws.Range("L86").Value = 25000
$ws.Range("N86").Value = -27372
$ws.Range("H89").Value = 25000
$ws.Range("J89").Value = 25000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86856
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1378
$ws.Range("N102").Value = -6244
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2387.8
$ws.Range("I20").Value = 1797.8462
$ws.Range("J20").Value = 3483.4285
$ws.Range("K20").Value = 1797.8462
$ws.Range("L20").Value = 3483.4285
$ws.Range("M20").Value = -1550.8462
$ws.Range("N20").Value = -3977.4285
$ws.Range("H82").Value = 15666
$ws.Range("J82").Value = 47833.332
$ws.Range("L82").Value = 47833.332
$ws.Range("N82").Value = -48599.332
$ws.Range("H85").Value = 15666
$ws.Range("J85").Value = 47833.332
$ws.Range("L85").Value = 47833.332
$ws.Range("N85").Value = -50485.332
$ws.Range("H94").Value = 1284.2307
$ws.Range("I94").Value = 556.625
$ws.Range("J94").Value = 2448.4
$ws.Range("K94").Value = 556.625
$ws.Range("L94").Value = 2448.4
$ws.Range("M94").Value = -105.625
$ws.Range("N94").Value = -3350.4
$ws.Range("H125").Value = 19000
$ws.Range("J125").Value = 19000
$ws.Range("L125").Value = 19000
$ws.Range("N125").Value = -28840

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1309.8948
$ws.Range("I22").Value = 799
$ws.Range("J22").Value = 1877.5555
$ws.Range("K22").Value = 799
$ws.Range("L22").Value = 1877.5555
$ws.Range("M22").Value = -449
$ws.Range("N22").Value = -2577.5555
$ws.Range("H31").Value = 3954.875
$ws.Range("I31").Value = 2058.8572
$ws.Range("J31").Value = 5429.5557
$ws.Range("K31").Value = 2058.8572
$ws.Range("L31").Value = 5429.5557
$ws.Range("M31").Value = -1763.8572
$ws.Range("N31").Value = -6019.5557
$ws.Range("H34").Value = 3954.875
$ws.Range("I34").Value = 2058.8572
$ws.Range("J34").Value = 5429.5557
$ws.Range("K34").Value = 2058.8572
$ws.Range("L34").Value = 5429.5557
$ws.Range("M34").Value = -1856.8572
$ws.Range("N34").Value = -5833.5557
$ws.Range("H107").Value = 2594.9
$ws.Range("J107").Value = 3224.5
$ws.Range("L107").Value = 3224.5
$ws.Range("N107").Value = -7064.5
$ws.Range("H132").Value = 2252.7896
$ws.Range("I132").Value = 1425.3125
$ws.Range("K132").Value = 4275.9375
$ws.Range("M132").Value = -1745.9375
$ws.Range("H141").Value = 73599.086
$ws.Range("J141").Value = 73744.45
$ws.Range("L141").Value = 73744.45
$ws.Range("N141").Value = -84104.45

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2302.5557
$ws.Range("I3").Value = 651.8
$ws.Range("J3").Value = 4366
$ws.Range("K3").Value = 1955.4
$ws.Range("L3").Value = 13098
$ws.Range("M3").Value = -1843.4
$ws.Range("N3").Value = -13322
$ws.Range("H55").Value = 3250
$ws.Range("J55").Value = 3250
$ws.Range("L55").Value = 9750
$ws.Range("N55").Value = -10104
$ws.Range("H114").Value = 2485.8333
$ws.Range("I114").Value = 452.85715
$ws.Range("J114").Value = 5332
$ws.Range("K114").Value = 1358.57145
$ws.Range("L114").Value = 15996
$ws.Range("M114").Value = 1895.42855
$ws.Range("N114").Value = -22504

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5611.3335
$ws.Range("I70").Value = 4499
$ws.Range("J70").Value = 5750.375
$ws.Range("K70").Value = 4499
$ws.Range("L70").Value = 5750.375
$ws.Range("M70").Value = -4229
$ws.Range("N70").Value = -6290.375
$ws.Range("H73").Value = 5611.3335
$ws.Range("I73").Value = 4499
$ws.Range("J73").Value = 5750.375
$ws.Range("K73").Value = 4499
$ws.Range("L73").Value = 5750.375
$ws.Range("M73").Value = -3563
$ws.Range("N73").Value = -7622.375
$ws.Range("H102").Value = 2079.9048
$ws.Range("I102").Value = 1811.125
$ws.Range("K102").Value = 1811.125
$ws.Range("M102").Value = -189.125
$ws.Range("H132").Value = 729432.5
$ws.Range("I132").Value = 965593.1
$ws.Range("K132").Value = 2896779.3
$ws.Range("M132").Value = -2894249.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1105.1666
$ws.Range("I22").Value = 823.25
$ws.Range("K22").Value = 823.25
$ws.Range("M22").Value = -528.25
$ws.Range("H27").Value = 1105.1666
$ws.Range("I27").Value = 823.25
$ws.Range("K27").Value = 823.25
$ws.Range("M27").Value = -716.25
$ws.Range("H68").Value = 3166.6667
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 4333.3335
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 4333.3335
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -5831.3335
$ws.Range("H69").Value = 58790.625
$ws.Range("J69").Value = 45065.2
$ws.Range("L69").Value = 45065.2
$ws.Range("N69").Value = -46687.2
$ws.Range("H71").Value = 3166.6667
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 4333.3335
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 21666.6675
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -29154.6675
$ws.Range("H72").Value = 58790.625
$ws.Range("J72").Value = 45065.2
$ws.Range("L72").Value = 135195.6
$ws.Range("N72").Value = -143307.6
$ws.Range("H82").Value = 2209.1667
$ws.Range("I82").Value = 2051
$ws.Range("K82").Value = 2051
$ws.Range("M82").Value = -1690
$ws.Range("H85").Value = 2209.1667
$ws.Range("I85").Value = 2051
$ws.Range("K85").Value = 2051
$ws.Range("M85").Value = -803
$ws.Range("H100").Value = 2909.9
$ws.Range("I100").Value = 1670
$ws.Range("J100").Value = 4149.8
$ws.Range("K100").Value = 1670
$ws.Range("L100").Value = 4149.8
$ws.Range("M100").Value = -1129
$ws.Range("N100").Value = -5231.8
$ws.Range("H109").Value = 59999
$ws.Range("J109").Value = 59999
$ws.Range("L109").Value = 59999
$ws.Range("N109").Value = -62773
$ws.Range("H132").Value = 4215.7075
$ws.Range("I132").Value = 4483.1333
$ws.Range("J132").Value = 4061.423
$ws.Range("K132").Value = 13449.3999
$ws.Range("L132").Value = 12184.269
$ws.Range("M132").Value = -10919.3999
$ws.Range("N132").Value = -17244.269
$ws.Range("H136").Value = 1642.0769
$ws.Range("I136").Value = 1654.7778
$ws.Range("J136").Value = 1613.5
$ws.Range("K136").Value = 4964.3334
$ws.Range("L136").Value = 4840.5
$ws.Range("M136").Value = -2414.3334
$ws.Range("N136").Value = -9940.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 67065
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H122").Value = 53590.438
$ws.Range("I122").Value = 68287.336
$ws.Range("K122").Value = 204862.008
$ws.Range("M122").Value = -202412.008
$ws.Range("H126").Value = 4359.077
$ws.Range("J126").Value = 5499.5
$ws.Range("L126").Value = 16498.5
$ws.Range("N126").Value = -21438.5
$ws.Range("H132").Value = 1366.2174
$ws.Range("I132").Value = 1141.8422
$ws.Range("K132").Value = 3425.5266
$ws.Range("M132").Value = -895.5266000000001
$ws.Range("H136").Value = 17923142
$ws.Range("I136").Value = 32681092
$ws.Range("J136").Value = 2772.4285
$ws.Range("K136").Value = 98043276
$ws.Range("L136").Value = 8317.2855
$ws.Range("M136").Value = -98040726
$ws.Range("N136").Value = -13417.2855
